# Actualiza datos del paciente en el formulario "HOJA DE INGRESO Y EGRESO"
# para agregar codigo de barras

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apellidos, nombres y no. de expediente clinico
$ws.Range("A6").Value = "MARTINEZ"
$ws.Range("C6").Value = "GOMEZ"
$ws.Range("E6").Value = "ERIKA"
$ws.Range("G6").Value = "YAJAIRA"
$ws.Range("I6").Value = "2017-7032/201773480"

# Direccion actual (calle, municipio, departamento, telefono)
$ws.Range("A8").Value = ""
$ws.Range("D8").Value = "ALDEA POXTE "
$ws.Range("F8").Value = "POPTUN"
$ws.Range("H8").Value = "PETEN"
$ws.Range("J8").NumberFormat = "@"
$ws.Range("J8").Value = "45325162"

# Fecha de nacimiento, edad, lugar de nacimiento, sexo
$ws.Range("A12").NumberFormat = "@"
$ws.Range("A12").Value = "1999-11-07"
$ws.Range("F12").NumberFormat = "@"
$ws.Range("F12").Value = "18"
$ws.Range("H12").Value = "POPTUN"
$ws.Range("J12").Value = "Femenino"

# Ocupacion, nacionalidad, no. de cedula
$ws.Range("D14").Value = "ESTUDIANTE"
$ws.Range("F14").Value = "GAUTEMALTECA"
$ws.Range("H14").NumberFormat = "@"
$ws.Range("H14").Value = "3306770471712"

# Nombre del padre y de la madre
$ws.Range("A18").Value = "EDIN ESTUARDO MARTINEZ"
$ws.Range("F18").Value = "ROSA MELIDA GOMEZ"

# En caso de emergencia notificar a: nombre, parentesco, direccion, telefono
$ws.Range("A20").Value = "ROSELIA VASQUEZ"
$ws.Range("F20").Value = "ABUELA"
$ws.Range("H20").Value = "POPTUN"
$ws.Range("J20").NumberFormat = "@"
$ws.Range("J20").Value = "48852148"

# Fecha de ingreso y hora
$ws.Range("A24").NumberFormat = "@"
$ws.Range("A24").Value = "20/11/2017"
$ws.Range("C24").NumberFormat = "@"
$ws.Range("C24").Value = "15:20:32"
